# SeiaOrg.xlsx update: add Ohio ("OH") row data pulled from
# https://seia.org/states-map, trim the now-superfluous New Jersey detail
# columns / "Projected Ranks" column, and drop the threaded comments that
# went with the removed columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the threaded comments attached to the headers of the columns
#     that are going away (D/F/G/H/I), keeping B1 (ITC), C1 (Solar
#     Rankings) and E1 (Percent electricity).
foreach ($ref in @("D1", "F1", "G1", "H1", "I1")) {
    $ws.Range($ref).CommentThreaded.Delete()
}

# --- New Jersey (row 2) no longer carries the detail columns C:J -- only
#     State (A) and ITC (B) remain populated for that row.
$ws.Range("C2:J2").ClearContents()

# --- Ohio (row 3) gains the full set of stats in C3:I3.
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = 95263
$ws.Range("E3").Value = 0.49
$ws.Range("F3").Value = 251
$ws.Range("G3").Value = 1300000000
$ws.Range("G3").NumberFormat = "#,##0"
$ws.Range("H3").Value = 0.36
$ws.Range("H3").NumberFormat = "0%"
$ws.Range("I3").Value = 5596

# --- The "Projected Ranks" column (J) is no longer used anywhere in the
#     sheet; clearing J1 drops the column from the used range / dimension
#     and prunes the now-unreferenced shared string.
$ws.Range("J1").ClearContents()

# --- Restore the cursor to where the author left it.
$ws.Range("C5").Select() | Out-Null
